# -----------------------------------------------------------------------
# Recipient Heat Fuel Fractions.xlsx - apply the "Adds district heat lever
# and updates BAU heat fractions" commit.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "About" sheet
# ------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Clear out the old "Notes:" block (rows 5-9) - it is being rebuilt two
# rows lower, with extra source-description rows inserted above it.
$about.Range("A5:B9").Clear() | Out-Null

# --- Source block -----------------------------------------------------
$about.Range("B3").Value = "Residential and Commercial"
$about.Range("B3").Font.Bold = $true
$about.Range("B3").Interior.Color = 12566463

$about.Range("B4").Value = "Based on FORECAST simulation, executed by TEP Energy"

$about.Range("B5").Value = "https://www.forecast-model.eu"
$about.Range("B5").HorizontalAlignment = -4131

$about.Range("B6").HorizontalAlignment = -4131

# --- Notes block (now starting two rows lower) -------------------------
$about.Range("A7").Value = "Notes:"
$about.Range("A7").Font.Bold = $true

$about.Range("A8").Value = "This variable specifies the recipient fuel for the Fraction of District Heat Fuel Use Shifted"
$about.Range("A9").Value = "to Other Fuels policy."

$about.Range("A10").Font.Bold = $true

$about.Range("A11").Value = "We assume a distribution of shifting to electricity via large scale electric heat pumps and"
$about.Range("A12").Value = "replacement of natural gas with hydrogen-burning facilities."

$about.Range("B17").Select() | Out-Null

# ------------------------------------------------------------------
# 2) "RHFF" sheet - update BAU heat fuel fractions
# ------------------------------------------------------------------
$rhff = $wb.Worksheets.Item("RHFF")

$rhff.Range("D2").Value = 0.9
$rhff.Range("E2").Value = 0
$rhff.Range("E5").Value = 1
$rhff.Range("D11").Value = 0.1

$rhff.Range("E6").Select() | Out-Null
$rhff.Activate() | Out-Null
